$d = $word.ActiveDocument

# Locate the paragraph that contains the exact sentence "Kassasjoner er registrert."
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Kassasjoner er registrert.`r") {
        $target = $p
    }
}

if ($target -eq $null) {
    Write-Host "Target paragraph not found!"
} else {
    $insPoint = $d.Range($target.Range.Start, $target.Range.Start)

    # Insert a new run "Varsel: " (with preserved trailing space) right before
    # the existing "Kassasjoner er registrert." run, as its own separate run
    # (matching the OOXML produced by the diff: two sibling <w:r> elements).
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Varsel: </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $inserted = $false
    try {
        $insPoint.InsertXML($xml)
        $inserted = $true
    } catch {
        $inserted = $false
    }

    if (-not $inserted) {
        # Fallback: plain text insertion (still yields correct text content
        # even if Word happens to merge it into the neighbouring run).
        $insPoint.InsertBefore("Varsel: ")
    }

    Write-Host "Updated paragraph text:" $target.Range.Text
}
